$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6989.3335
$ws.Range("J112").Value = 7313.8
$ws.Range("L112").Value = 21941.4
$ws.Range("N112").Value = -24157.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1286.7391
$ws.Range("I129").Value = 487.85715
$ws.Range("J129").Value = 1636.25
$ws.Range("K129").Value = 1463.57145
$ws.Range("L129").Value = 4908.75
$ws.Range("M129").Value = 3536.42855
$ws.Range("N129").Value = -14908.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1341.8077
$ws.Range("I137").Value = 1257.7646
$ws.Range("J137").Value = 1500.5555
$ws.Range("K137").Value = 3773.2938
$ws.Range("L137").Value = 4501.666499999999
$ws.Range("M137").Value = -1223.2938
$ws.Range("N137").Value = -9601.666499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11536.56
$ws.Range("I32").Value = 11768.23
$ws.Range("K32").Value = 11768.23
$ws.Range("M32").Value = -11481.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3012.0356
$ws.Range("I61").Value = 2268.3
$ws.Range("J61").Value = 4871.375
$ws.Range("K61").Value = 2268.3
$ws.Range("L61").Value = 4871.375
$ws.Range("M61").Value = -2056.3
$ws.Range("N61").Value = -5295.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1176.68
$ws.Range("I74").Value = 1075.5625
$ws.Range("J74").Value = 1356.4445
$ws.Range("K74").Value = 1075.5625
$ws.Range("L74").Value = 1356.4445
$ws.Range("M74").Value = -201.5625
$ws.Range("N74").Value = -3104.4445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1176.68
$ws.Range("I77").Value = 1075.5625
$ws.Range("J77").Value = 1356.4445
$ws.Range("K77").Value = 5377.8125
$ws.Range("L77").Value = 6782.2225
$ws.Range("M77").Value = -1009.8125
$ws.Range("N77").Value = -15518.2225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 251974.75
$ws.Range("I102").Value = 2449.5
$ws.Range("K102").Value = 2449.5
$ws.Range("M102").Value = -827.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1340.2
$ws.Range("I122").Value = 1425.25
$ws.Range("K122").Value = 4275.75
$ws.Range("M122").Value = -1825.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3012.0356
$ws.Range("I136").Value = 2268.3
$ws.Range("J136").Value = 4871.375
$ws.Range("K136").Value = 6804.900000000001
$ws.Range("L136").Value = 14614.125
$ws.Range("M136").Value = -4254.900000000001
$ws.Range("N136").Value = -19714.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2041.7241
$ws.Range("I31").Value = 1276.6666
$ws.Range("J31").Value = 4050
$ws.Range("K31").Value = 1276.6666
$ws.Range("L31").Value = 4050
$ws.Range("M31").Value = -981.6666
$ws.Range("N31").Value = -4640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2041.7241
$ws.Range("I34").Value = 1276.6666
$ws.Range("J34").Value = 4050
$ws.Range("K34").Value = 1276.6666
$ws.Range("L34").Value = 4050
$ws.Range("M34").Value = -1074.6666
$ws.Range("N34").Value = -4454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2021.4445
$ws.Range("J94").Value = 2230.5833
$ws.Range("L94").Value = 2230.5833
$ws.Range("N94").Value = -3132.5833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 271636.38
$ws.Range("I132").Value = 387165.62
$ws.Range("J132").Value = 2068.0667
$ws.Range("K132").Value = 1161496.86
$ws.Range("L132").Value = 6204.2001
$ws.Range("M132").Value = -1158966.86
$ws.Range("N132").Value = -11264.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 706.7646999999999
$ws.Range("I2").Value = 1091.5
$ws.Range("J2").Value = 157.14285
$ws.Range("K2").Value = 6549
$ws.Range("L2").Value = 942.8571000000001
$ws.Range("M2").Value = -6436
$ws.Range("N2").Value = -1168.8571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 7287.375
$ws.Range("J35").Value = 8285.571
$ws.Range("L35").Value = 24856.713
$ws.Range("N35").Value = -25432.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 980.8182
$ws.Range("I68").Value = 1020
$ws.Range("J68").Value = 958.4286
$ws.Range("K68").Value = 3060
$ws.Range("L68").Value = 2875.2858
$ws.Range("M68").Value = -2249
$ws.Range("N68").Value = -4497.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 980.8182
$ws.Range("I71").Value = 1020
$ws.Range("J71").Value = 958.4286
$ws.Range("K71").Value = 9180
$ws.Range("L71").Value = 8625.857399999999
$ws.Range("M71").Value = -5124
$ws.Range("N71").Value = -16737.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 453.9091
$ws.Range("J92").Value = 393
$ws.Range("L92").Value = 1179
$ws.Range("N92").Value = -3675

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1746.3334
$ws.Range("J98").Value = 761.6
$ws.Range("L98").Value = 2284.8
$ws.Range("N98").Value = -5280.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 4933.3335
$ws.Range("I99").Value = 2400
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 7200
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -4954
$ws.Range("N99").Value = -34492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4702
$ws.Range("I116").Value = 904
$ws.Range("J116").Value = 8500
$ws.Range("K116").Value = 2712
$ws.Range("L116").Value = 25500
$ws.Range("M116").Value = 730
$ws.Range("N116").Value = -32384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1395.8
$ws.Range("I122").Value = 996.6667
$ws.Range("J122").Value = 1994.5
$ws.Range("K122").Value = 8970.0003
$ws.Range("L122").Value = 17950.5
$ws.Range("M122").Value = -6520.0003
$ws.Range("N122").Value = -22850.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 56255224
$ws.Range("I122").Value = 62505500
$ws.Range("J122").Value = 50004950
$ws.Range("K122").Value = 187516500
$ws.Range("L122").Value = 150014850
$ws.Range("M122").Value = -187514050
$ws.Range("N122").Value = -150019750

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4695.75
$ws.Range("I132").Value = 4077.8462
$ws.Range("J132").Value = 5843.2856
$ws.Range("K132").Value = 12233.5386
$ws.Range("L132").Value = 17529.8568
$ws.Range("M132").Value = -9703.5386
$ws.Range("N132").Value = -22589.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2990.2646
$ws.Range("I136").Value = 2260.7917
$ws.Range("J136").Value = 4741
$ws.Range("K136").Value = 6782.375100000001
$ws.Range("L136").Value = 14223
$ws.Range("M136").Value = -4232.375100000001
$ws.Range("N136").Value = -19323

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10418925
$ws.Range("I122").Value = 14707185
$ws.Range("J122").Value = 4811201
$ws.Range("K122").Value = 44121555
$ws.Range("L122").Value = 14433603
$ws.Range("M122").Value = -44119105
$ws.Range("N122").Value = -14438503

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 35928.75
$ws.Range("J123").Value = 35928.75
$ws.Range("L123").Value = 35928.75
$ws.Range("N123").Value = -45728.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1873.8077
$ws.Range("I132").Value = 1738.5
$ws.Range("J132").Value = 2090.3
$ws.Range("K132").Value = 5215.5
$ws.Range("L132").Value = 6270.900000000001
$ws.Range("M132").Value = -2685.5
$ws.Range("N132").Value = -11330.9
